# Added diagrams for AddOrderCommand and edited some descriptions
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Reposition/resize the "Elbow Connector 63" (id=53) bent connector ---
$connector = $s.Shapes.Item(19)
$connector.Left = 50.82700987401575
$connector.Top = 325.61558555118114
$connector.Width = 271.0236360472441
$connector.Height = 43.09630021259842

# --- Merge split text runs into single runs with updated wording ---
# "Order" + "ListPanel" -> "OrderListPanel" (Rectangle 11, id=65, inside Group 14)
$group14 = $s.Shapes.Item(41)
$orderListPanel = $group14.GroupItems.Item(1)
$orderListPanel.TextFrame.TextRange.Text = "x"
$orderListPanel.TextFrame.TextRange.Text = "OrderListPanel"

# "Order" + "Card" -> "OrderCard" (Rectangle 11, id=66, inside Group 14)
$orderCard = $group14.GroupItems.Item(2)
$orderCard.TextFrame.TextRange.Text = "x"
$orderCard.TextFrame.TextRange.Text = "OrderCard"

# "Right" + "Panel" -> "RightPanel" (Rectangle 11, id=98)
$rightPanel = $s.Shapes.Item(45)
$rightPanel.TextFrame.TextRange.Text = "x"
$rightPanel.TextFrame.TextRange.Text = "RightPanel"

# "Person" + "Panel" -> "PersonPanel" (Rectangle 11, id=123)
$personPanel = $s.Shapes.Item(49)
$personPanel.TextFrame.TextRange.Text = "x"
$personPanel.TextFrame.TextRange.Text = "PersonPanel"

# --- Remove the extra (empty) slide 2 ---
$p.Slides.Item(2).Delete()
